$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (old B -> C, old C -> D), matching
# the formatting (width) of column A, and shift the header/query columns over.
$ws.Columns("B:B").Insert()
$ws.Columns("B:B").ColumnWidth = $ws.Columns("A:A").ColumnWidth

# New header for the inserted column.
$ws.Range("B1").Value = "StatQuery"

# New stats query for the inserted column, styled the same as A2 (wrap text).
$ws.Range("B2").Value = "MATCH (t:clinical_trial)<--(a:arm)<--(c:case)<--(s:specimen)<--(:assignment_report) WITH DISTINCT c AS c, t ,a, s WHERE c.race IN ['WHITE'] OPTIONAL MATCH (s)<-[*]-(f:file) RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(t.clinical_trial_designation)) as number_of_trial"
$ws.Range("B2").WrapText = $true

# Match the final active selection recorded for the sheet.
$ws.Range("A4").Select() | Out-Null
